# Update the division-problem answers in the table to match the new
# generated output set. Each old value is unique in the document, so a
# simple Find/Replace per pair is safe and order-independent.

$d = $word.ActiveDocument

$pairs = @(
    @("98÷7=14, 0", "75÷4=18, 3"),
    @("76÷9=8, 4", "39÷8=4, 7"),
    @("32÷2=16, 0", "70÷9=7, 7"),
    @("30÷3=10, 0", "99÷5=19, 4"),
    @("78÷2=39, 0", "70÷5=14, 0"),
    @("71÷5=14, 1", "83÷5=16, 3"),
    @("40÷3=13, 1", "11÷5=2, 1"),
    @("35÷7=5, 0", "13÷4=3, 1"),
    @("26÷8=3, 2", "40÷7=5, 5"),
    @("83÷9=9, 2", "37÷7=5, 2"),
    @("97÷8=12, 1", "73÷7=10, 3"),
    @("67÷9=7, 4", "44÷8=5, 4"),
    @("99÷9=11, 0", "65÷9=7, 2"),
    @("25÷4=6, 1", "39÷5=7, 4"),
    @("58÷6=9, 4", "73÷4=18, 1"),
    @("41÷8=5, 1", "70÷8=8, 6"),
    @("49÷5=9, 4", "58÷2=29, 0"),
    @("94÷5=18, 4", "60÷4=15, 0"),
    @("53÷2=26, 1", "67÷3=22, 1"),
    @("74÷8=9, 2", "76÷3=25, 1"),
    @("21÷8=2, 5", "17÷6=2, 5"),
    @("51÷6=8, 3", "74÷2=37, 0"),
    @("12÷6=2, 0", "74÷2=37, 0"),
    @("53÷3=17, 2", "11÷5=2, 1"),
    @("87÷8=10, 7", "89÷3=29, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
